# Fixed 3D model issues
# The part number for "Switch SP4T" was corrected from SKY13380-350LF to PE42440.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = "PE42440"

# Update the active cell selection to match the saved state in the workbook.
$ws.Range("E12").Select()
